# "Penalty Reward System" edit (unfinished, per the author's commit message):
#  - Forecast Comparison sheet: shift each week's Week_Start_Date forward by
#    one week, and zero out the MyForecast column for every data row.
#  - Summary sheet: update a handful of metrics to new values.
#
# All of the target cells are stored as plain text (not real numbers/dates)
# in the source workbook, even though several look like numbers or dates.
# A leading apostrophe forces Excel to keep them as text instead of
# auto-converting to a date serial / numeric value, matching the original
# inlineStr-as-text representation.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Forecast Comparison" ----
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = "'" + $newDates[$i]
    $ws1.Cells.Item($row, 4).Value = 0
}

# ---- Sheet 2: "Summary" ----
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-05"
$ws2.Range("B4").Value = "'69"
$ws2.Range("B8").Value = "1554 units"
$ws2.Range("B9").Value = "'5"
$ws2.Range("B10").Value = "'3"
$ws2.Range("B11").Value = "'2"
$ws2.Range("B12").Value = "'0"
$ws2.Range("B14").Value = "'0"
$ws2.Range("B15").Value = "'2025-02-09"
